# Review with product owner
# Append a new journal entry (row 51) documenting the PO review, matching
# the style/shape of the previous entries, then extend the table/autofilter
# to cover the new row and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row of data (row 51) ------------------------------------
# Clone the date formatting from A50 (builtin date number format) onto A51
# before writing the new date value into it.
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A51").PasteSpecial(-4122) | Out-Null

$ws.Range("A51").Value = 44694
$ws.Range("B51").Value = "Review"
$ws.Range("C51").Value = 0.5
$ws.Range("D51").Value = "Sprint review avec chef de projet "
$ws.Range("E51").Value = "Quelques trucs a réviser. Voir dans documentation/review"

# Match the taller row height used by similar two-line entries.
$ws.Rows.Item(51).RowHeight = 30

# --- Expand the table / autofilter to include the new row ----------------
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Resize($ws.Range("A1:F51")) | Out-Null

# --- Update the view: scroll position and active selection ---------------
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("A52").Select() | Out-Null
